{"js": "// Insert a new first paragraph \"Team 8\" (yellow-highlighted) above the\n// existing first paragraph of the document body, matching the\n// \"Reformat & include team name\" commit.\nconst body = context.document.body;\n\n// Insert a brand-new paragraph before everything currently in the body.\nconst newPara = body.insertParagraph(\"Team 8\", \"Start\");\n\n// Highlight both the paragraph mark (pPr/rPr) and the run text itself,\n// same as the rest of the \"Names\" block above it.\nnewPara.font.highlightColor = \"Yellow\";\n\nawait context.sync();\n", "ps1": "# Insert a new first paragraph \"Team 8\" (yellow-highlighted) above the\n# existing first paragraph of the document body, matching the\n# \"Reformat & include team name\" commit.\n$d = $word.ActiveDocument\n\n# Create a new empty paragraph right before the current first paragraph.\n$d.Paragraphs(1).Range.InsertParagraphBefore()\n\n# That new paragraph is now Paragraphs(1); fill it in and highlight it\n# the same way the rest of the \"Names\" block above it is highlighted.\n$newPara = $d.Paragraphs(1)\n$newPara.Range.Text = \"Team 8\"\n$newPara.Range.HighlightColorIndex = \"Yellow\"\n"}
